$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text/content edits (TFL -> AF variant row, 2nd data row becomes "AF")
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "LogIn Alis TFL"
$ws.Range("B3").Value = "LogIn Alis AF"
$ws.Range("E3").Value = "a"
$ws.Range("F3").Value = "af_7000_defs"

# ---------------------------------------------------------------------------
# 2. Turn C3 into a real (external) hyperlink pointing at the new AF env,
#    while giving the whole C2:C7 "URL" column the same Hyperlink look
#    (underlined, theme Hyperlink colour) that Excel applies when you use
#    Insert > Hyperlink.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "http://alis-deploy15:8094/tfl_pl_env8/alis#alis", "", "http://alis-deploy15:8094/tfl_pl_env8/alis#alis")
$ws.Hyperlinks.Add($ws.Range("C4"), "http://alis-deploy15:8094/tfl_pl_env8/alis#alis", "", "http://alis-deploy15:8094/tfl_pl_env8/alis#alis")
$ws.Hyperlinks.Add($ws.Range("C5"), "http://alis-deploy15:8094/tfl_pl_env8/alis#alis", "", "http://alis-deploy15:8094/tfl_pl_env8/alis#alis")
$ws.Hyperlinks.Add($ws.Range("C6"), "http://alis-deploy15:8094/tfl_pl_env8/alis#alis", "", "http://alis-deploy15:8094/tfl_pl_env8/alis#alis")
$ws.Hyperlinks.Add($ws.Range("C7"), "http://alis-deploy15:8094/tfl_pl_env8/alis#alis", "", "http://alis-deploy15:8094/tfl_pl_env8/alis#alis")
$ws.Hyperlinks.Add($ws.Range("C3"), "http://alis-alf-app01:8080/af_pl_env1/alis1#alis", "alis", "", "http://alis-alf-app01:8080/af_pl_env1/alis1 - alis")
$ws.Range("C3").Value = "http://alis-alf-app01:8080/af_pl_env1/alis1#alis"

# Only C3 should stay a "live" hyperlink entry - drop the helper links we
# added to C2/C4:C7 again, this leaves their Hyperlink-style formatting
# (font/underline/colour) intact while un-registering the jump target.
for ($i = $ws.Hyperlinks.Count; $i -ge 1; $i--) {
    $h = $ws.Hyperlinks.Item($i)
    if ($h.Range.Address() -ne '$C$3') {
        $h.Delete()
    }
}

$ws.Range("C2:C7").HorizontalAlignment = -4131
$ws.Range("C2:C7").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. F3 ("af_7000_defs") gets its own small Verdana/centered look.
# ---------------------------------------------------------------------------
$ws.Range("F3").Font.Size = 8
$ws.Range("F3").Font.Name = "Verdana"
$ws.Range("F3").HorizontalAlignment = -4108
$ws.Range("F3").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Column C got a touch narrower, selection moved to B4.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 44
$ws.Range("B4").Select()
